$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1) - give it the value first, then copy the
# formatting from an existing header cell (G1) so it picks up the same
# shared cell style (bold font, borders, centered alignment) rather than
# minting a brand new (duplicate) style.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New data value for the "Save" column (H2) - plain number, unstyled like
# the other numeric cells in the row.
$ws.Range("H2").Value = 0
